$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 4997
$ws.Cells.Item(86, 10).Value = 4997
$ws.Cells.Item(86, 12).Value = 4997
$ws.Cells.Item(86, 14).Value = -7243
$ws.Cells.Item(89, 8).Value = 4997
$ws.Cells.Item(89, 10).Value = 4997
$ws.Cells.Item(89, 12).Value = 24985
$ws.Cells.Item(89, 14).Value = -36217
$ws.Cells.Item(99, 8).Value = 739.4
$ws.Cells.Item(99, 9).Value = 624.25
$ws.Cells.Item(99, 11).Value = 1872.75
$ws.Cells.Item(99, 13).Value = -374.75
$ws.Cells.Item(132, 8).Value = 245317.25
$ws.Cells.Item(132, 9).Value = 1270.5405
$ws.Cells.Item(132, 11).Value = 3811.6215
$ws.Cells.Item(132, 13).Value = -1281.6215
$ws.Cells.Item(137, 8).Value = 5156.1665
$ws.Cells.Item(137, 9).Value = 3235.1667
$ws.Cells.Item(137, 10).Value = 7077.1665
$ws.Cells.Item(137, 11).Value = 9705.500100000001
$ws.Cells.Item(137, 12).Value = 21231.4995
$ws.Cells.Item(137, 13).Value = -7155.500100000001
$ws.Cells.Item(137, 14).Value = -26331.4995
$ws.Cells.Item(138, 8).Value = 4285.023
$ws.Cells.Item(138, 9).Value = 2476.25
$ws.Cells.Item(138, 10).Value = 4686.972
$ws.Cells.Item(138, 11).Value = 7428.75
$ws.Cells.Item(138, 12).Value = 14060.916
$ws.Cells.Item(138, 13).Value = -2288.75
$ws.Cells.Item(138, 14).Value = -24340.916
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1360.7778
$ws.Cells.Item(2, 9).Value = 1371.1143
$ws.Cells.Item(2, 11).Value = 1371.1143
$ws.Cells.Item(2, 13).Value = -1258.1143
$ws.Cells.Item(45, 8).Value = 4810.5713
$ws.Cells.Item(45, 9).Value = 4332.2
$ws.Cells.Item(45, 11).Value = 4332.2
$ws.Cells.Item(45, 13).Value = -3955.2
$ws.Cells.Item(61, 8).Value = 4061.516
$ws.Cells.Item(61, 9).Value = 3852.1667
$ws.Cells.Item(61, 11).Value = 3852.1667
$ws.Cells.Item(61, 13).Value = -3640.1667
$ws.Cells.Item(74, 8).Value = 793.8125
$ws.Cells.Item(74, 9).Value = 750.3077
$ws.Cells.Item(74, 11).Value = 750.3077
$ws.Cells.Item(74, 13).Value = 123.6923
$ws.Cells.Item(77, 8).Value = 793.8125
$ws.Cells.Item(77, 9).Value = 750.3077
$ws.Cells.Item(77, 11).Value = 3751.5385
$ws.Cells.Item(77, 13).Value = 616.4615000000003
$ws.Cells.Item(116, 8).Value = 1360.7778
$ws.Cells.Item(116, 9).Value = 1371.1143
$ws.Cells.Item(116, 11).Value = 1371.1143
$ws.Cells.Item(116, 13).Value = 922.8857
$ws.Cells.Item(122, 8).Value = 2634.775
$ws.Cells.Item(122, 9).Value = 2561.5588
$ws.Cells.Item(122, 11).Value = 7684.676399999999
$ws.Cells.Item(122, 13).Value = -5234.676399999999
$ws.Cells.Item(136, 8).Value = 4061.516
$ws.Cells.Item(136, 9).Value = 3852.1667
$ws.Cells.Item(136, 11).Value = 11556.5001
$ws.Cells.Item(136, 13).Value = -9006.500100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1360.7778
$ws.Cells.Item(3, 9).Value = 1371.1143
$ws.Cells.Item(3, 11).Value = 1371.1143
$ws.Cells.Item(3, 13).Value = -1257.1143
$ws.Cells.Item(42, 8).Value = 175945
$ws.Cells.Item(42, 10).Value = 175945
$ws.Cells.Item(42, 12).Value = 175945
$ws.Cells.Item(42, 14).Value = -176601
$ws.Cells.Item(76, 8).Value = 17746.5
$ws.Cells.Item(76, 10).Value = 17746.5
$ws.Cells.Item(76, 12).Value = 17746.5
$ws.Cells.Item(76, 14).Value = -18376.5
$ws.Cells.Item(79, 8).Value = 17746.5
$ws.Cells.Item(79, 10).Value = 17746.5
$ws.Cells.Item(79, 12).Value = 17746.5
$ws.Cells.Item(79, 14).Value = -19930.5
$ws.Cells.Item(86, 8).Value = 2824.7144
$ws.Cells.Item(86, 9).Value = 2824.7144
$ws.Cells.Item(86, 11).Value = 2824.7144
$ws.Cells.Item(86, 13).Value = -1701.7144
$ws.Cells.Item(89, 8).Value = 2824.7144
$ws.Cells.Item(89, 9).Value = 2824.7144
$ws.Cells.Item(89, 11).Value = 14123.572
$ws.Cells.Item(89, 13).Value = -8507.572
$ws.Cells.Item(105, 8).Value = 3361.3333
$ws.Cells.Item(105, 9).Value = 1769.7142
$ws.Cells.Item(105, 10).Value = 4016.7058
$ws.Cells.Item(105, 11).Value = 1769.7142
$ws.Cells.Item(105, 12).Value = 4016.7058
$ws.Cells.Item(105, 13).Value = -22.71419999999989
$ws.Cells.Item(105, 14).Value = -7510.7058
$ws.Cells.Item(107, 8).Value = 664
$ws.Cells.Item(107, 9).Value = 645.3889
$ws.Cells.Item(107, 11).Value = 645.3889
$ws.Cells.Item(107, 13).Value = 1274.6111
$ws.Cells.Item(134, 8).Value = 4323.6665
$ws.Cells.Item(134, 10).Value = 4323.6665
$ws.Cells.Item(134, 12).Value = 12970.9995
$ws.Cells.Item(134, 14).Value = -18040.9995
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 11434529
$ws.Cells.Item(6, 9).Value = 16000140
$ws.Cells.Item(6, 11).Value = 16000140
$ws.Cells.Item(6, 13).Value = -16000027
$ws.Cells.Item(7, 8).Value = 91445.82000000001
$ws.Cells.Item(7, 9).Value = 111155.78
$ws.Cells.Item(7, 11).Value = 111155.78
$ws.Cells.Item(7, 13).Value = -111042.78
$ws.Cells.Item(17, 8).Value = 9999.5
$ws.Cells.Item(17, 9).Value = 9999
$ws.Cells.Item(17, 11).Value = 9999
$ws.Cells.Item(17, 13).Value = -9825
$ws.Cells.Item(31, 8).Value = 3622.75
$ws.Cells.Item(31, 9).Value = 2326.1035
$ws.Cells.Item(31, 10).Value = 4586.923
$ws.Cells.Item(31, 11).Value = 2326.1035
$ws.Cells.Item(31, 12).Value = 4586.923
$ws.Cells.Item(31, 13).Value = -2031.1035
$ws.Cells.Item(31, 14).Value = -5176.923
$ws.Cells.Item(34, 8).Value = 3622.75
$ws.Cells.Item(34, 9).Value = 2326.1035
$ws.Cells.Item(34, 10).Value = 4586.923
$ws.Cells.Item(34, 11).Value = 2326.1035
$ws.Cells.Item(34, 12).Value = 4586.923
$ws.Cells.Item(34, 13).Value = -2124.1035
$ws.Cells.Item(34, 14).Value = -4990.923
$ws.Cells.Item(51, 8).Value = 50000
$ws.Cells.Item(51, 10).Value = 50000
$ws.Cells.Item(51, 12).Value = 50000
$ws.Cells.Item(51, 14).Value = -51472
$ws.Cells.Item(58, 8).Value = 1113.8334
$ws.Cells.Item(58, 9).Value = 1136.6
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 1136.6
$ws.Cells.Item(58, 12).Value = 1000
$ws.Cells.Item(58, 13).Value = -933.5999999999999
$ws.Cells.Item(58, 14).Value = -1406
$ws.Cells.Item(59, 8).Value = 45039.824
$ws.Cells.Item(59, 9).Value = 26000
$ws.Cells.Item(59, 10).Value = 47578.465
$ws.Cells.Item(59, 11).Value = 26000
$ws.Cells.Item(59, 12).Value = 47578.465
$ws.Cells.Item(59, 13).Value = -24855
$ws.Cells.Item(59, 14).Value = -49868.465
$ws.Cells.Item(60, 8).Value = 7296.5
$ws.Cells.Item(61, 8).Value = 50000
$ws.Cells.Item(61, 10).Value = 50000
$ws.Cells.Item(61, 12).Value = 50000
$ws.Cells.Item(61, 14).Value = -50696
$ws.Cells.Item(105, 8).Value = 1877.3529
$ws.Cells.Item(105, 9).Value = 1760.6666
$ws.Cells.Item(105, 10).Value = 2157.4
$ws.Cells.Item(105, 11).Value = 1760.6666
$ws.Cells.Item(105, 12).Value = 2157.4
$ws.Cells.Item(105, 13).Value = -13.66660000000002
$ws.Cells.Item(105, 14).Value = -5651.4
$ws.Cells.Item(136, 8).Value = 1113.8334
$ws.Cells.Item(136, 9).Value = 1136.6
$ws.Cells.Item(136, 10).Value = 1000
$ws.Cells.Item(136, 11).Value = 3409.8
$ws.Cells.Item(136, 12).Value = 3000
$ws.Cells.Item(136, 13).Value = -859.7999999999997
$ws.Cells.Item(136, 14).Value = -8100
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 116.666664
$ws.Cells.Item(15, 9).Value = 16
$ws.Cells.Item(15, 11).Value = 48
$ws.Cells.Item(15, 13).Value = 92
$ws.Cells.Item(56, 8).Value = 5483
$ws.Cells.Item(56, 9).Value = 5483
$ws.Cells.Item(56, 11).Value = 5483
$ws.Cells.Item(56, 13).Value = -4953
$ws.Cells.Item(59, 8).Value = 566.3333
$ws.Cells.Item(59, 9).Value = 566.3333
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 11).Value = 1698.9999
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()
$ws.Cells.Item(59, 13).Value = -1158.9999
$ws.Cells.Item(75, 8).Value = 4001.5
$ws.Cells.Item(75, 10).Value = 3990
$ws.Cells.Item(75, 12).Value = 11970
$ws.Cells.Item(75, 14).Value = -13966
$ws.Cells.Item(78, 8).Value = 4001.5
$ws.Cells.Item(78, 10).Value = 3990
$ws.Cells.Item(78, 12).Value = 35910
$ws.Cells.Item(78, 14).Value = -45894
$ws.Cells.Item(103, 8).Value = 483.66666
$ws.Cells.Item(103, 9).Value = 483.66666
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 1450.99998
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).ClearContents()
$ws.Cells.Item(103, 13).Value = -571.9999800000001
$ws.Cells.Item(132, 8).Value = 3724
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 12).Value = 36000
$ws.Cells.Item(132, 14).Value = -41060
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2853.5715
$ws.Cells.Item(122, 10).Value = 2994
$ws.Cells.Item(122, 12).Value = 8982
$ws.Cells.Item(122, 14).Value = -13882
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3935.6667
$ws.Cells.Item(7, 9).Value = 3917.8572
$ws.Cells.Item(7, 11).Value = 3917.8572
$ws.Cells.Item(7, 13).Value = -3805.8572
$ws.Cells.Item(25, 8).Value = 7159291.5
$ws.Cells.Item(40, 8).Value = 3632.625
$ws.Cells.Item(40, 9).Value = 3541.5334
$ws.Cells.Item(40, 11).Value = 3541.5334
$ws.Cells.Item(40, 13).Value = -3405.5334
$ws.Cells.Item(122, 8).Value = 3187.8125
$ws.Cells.Item(122, 9).Value = 2925.9167
$ws.Cells.Item(122, 11).Value = 8777.750100000001
$ws.Cells.Item(122, 13).Value = -6327.750100000001
$ws.Cells.Item(126, 8).Value = 3935.6667
$ws.Cells.Item(126, 9).Value = 3917.8572
$ws.Cells.Item(126, 11).Value = 11753.5716
$ws.Cells.Item(126, 13).Value = -9283.571599999999
$ws.Cells.Item(136, 8).Value = 2249.75
$ws.Cells.Item(136, 9).Value = 2249.75
$ws.Cells.Item(136, 11).Value = 6749.25
$ws.Cells.Item(136, 13).Value = -4199.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3667.4443
$ws.Cells.Item(96, 9).Value = 4000
$ws.Cells.Item(96, 11).Value = 4000
$ws.Cells.Item(96, 13).Value = -2627
$ws.Cells.Item(122, 8).Value = 3635.875
$ws.Cells.Item(122, 9).Value = 3460.2058
$ws.Cells.Item(122, 11).Value = 10380.6174
$ws.Cells.Item(122, 13).Value = -7930.617400000001
$ws.Cells.Item(124, 8).Value = 29666.666
$ws.Cells.Item(124, 10).Value = 29666.666
$ws.Cells.Item(124, 12).Value = 29666.666
$ws.Cells.Item(124, 14).Value = -39486.666
$ws.Cells.Item(126, 8).Value = 4547.9287
$ws.Cells.Item(126, 10).Value = 6926
$ws.Cells.Item(126, 12).Value = 20778
$ws.Cells.Item(126, 14).Value = -25718
